# Inserts two new price records (rows 539 and 540) into the daily pricing
# log, shifting all subsequent rows down by two positions.
#
# Row 539 (new): Fecha 44711 (2022-05-30), Primera, Volumen 162, Precio 18000,
#                Unidad "$/caja 12 unidades", Precio $/Kg 1500, Kg/unidad 12
# Row 540 (new): Fecha 44711 (2022-05-30), Segunda, Volumen 108, Precio 18000,
#                Unidad "$/caja 14 unidades", Precio $/Kg 1286, Kg/unidad 14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 539; this pushes the existing
# rows 539..651 down to 541..653 and keeps all of their data intact.
$ws.Rows("539:540").Insert()

# --- Row 539 -------------------------------------------------------------
$ws.Cells.Item(539, 1).Value = 3
$ws.Cells.Item(539, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(539, 3).Value = "Coquimbo"
$ws.Cells.Item(539, 4).Value = 44711
$ws.Cells.Item(539, 5).Value = 5
$ws.Cells.Item(539, 6).Value = "Fruta"
$ws.Cells.Item(539, 7).Value = 100108
$ws.Cells.Item(539, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(539, 9).Value = 100108005
$ws.Cells.Item(539, 10).Value = "Piña"
$ws.Cells.Item(539, 11).Value = "Caramelo"
$ws.Cells.Item(539, 12).Value = "Primera"
$ws.Cells.Item(539, 13).Value = 162
$ws.Cells.Item(539, 14).Value = 18000
$ws.Cells.Item(539, 15).Value = 18000
$ws.Cells.Item(539, 16).Value = 18000
$ws.Cells.Item(539, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(539, 18).Value = "Ecuador"
$ws.Cells.Item(539, 19).Value = 1500
$ws.Cells.Item(539, 20).Value = 12

# --- Row 540 -------------------------------------------------------------
$ws.Cells.Item(540, 1).Value = 3
$ws.Cells.Item(540, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(540, 3).Value = "Coquimbo"
$ws.Cells.Item(540, 4).Value = 44711
$ws.Cells.Item(540, 5).Value = 5
$ws.Cells.Item(540, 6).Value = "Fruta"
$ws.Cells.Item(540, 7).Value = 100108
$ws.Cells.Item(540, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(540, 9).Value = 100108005
$ws.Cells.Item(540, 10).Value = "Piña"
$ws.Cells.Item(540, 11).Value = "Caramelo"
$ws.Cells.Item(540, 12).Value = "Segunda"
$ws.Cells.Item(540, 13).Value = 108
$ws.Cells.Item(540, 14).Value = 18000
$ws.Cells.Item(540, 15).Value = 18000
$ws.Cells.Item(540, 16).Value = 18000
$ws.Cells.Item(540, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(540, 18).Value = "Ecuador"
$ws.Cells.Item(540, 19).Value = 1286
$ws.Cells.Item(540, 20).Value = 14
